$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 884.6667
$ws.Range("I53").Value = 1310.8889
$ws.Range("J53").Value = 565
$ws.Range("K53").Value = 1310.8889
$ws.Range("L53").Value = 565
$ws.Range("M53").Value = -673.8888999999999
$ws.Range("N53").Value = -1839

$ws.Range("H112").Value = 1571.6666
$ws.Range("J112").Value = 1671.25
$ws.Range("L112").Value = 5013.75
$ws.Range("N112").Value = -7229.75

$ws.Range("H116").Value = 2433.3333
$ws.Range("I116").Value = 2433.3333
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2433.3333
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1008.6667
$ws.Range("N116").ClearContents()

$ws.Range("H137").Value = 1180.8182
$ws.Range("I137").Value = 1198.9
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 3596.7
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -1046.7
$ws.Range("N137").Value = -8100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1800
$ws.Range("I2").Value = 1666.6666
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 1666.6666
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -1553.6666
$ws.Range("N2").Value = -2426

$ws.Range("H32").Value = 4722.25
$ws.Range("I32").Value = 4307.737
$ws.Range("J32").Value = 12598
$ws.Range("K32").Value = 4307.737
$ws.Range("L32").Value = 12598
$ws.Range("M32").Value = -4020.737
$ws.Range("N32").Value = -13172

$ws.Range("H35").Value = 686.75
$ws.Range("I35").Value = 686.75
$ws.Range("K35").Value = 686.75
$ws.Range("M35").Value = -280.75

$ws.Range("H45").Value = 1421947.1
$ws.Range("I45").Value = 1819599.4
$ws.Range("J45").Value = 1760.5714
$ws.Range("K45").Value = 1819599.4
$ws.Range("L45").Value = 1760.5714
$ws.Range("M45").Value = -1819222.4
$ws.Range("N45").Value = -2514.5714

$ws.Range("H61").Value = 3181.88
$ws.Range("I61").Value = 3231.125
$ws.Range("K61").Value = 3231.125
$ws.Range("M61").Value = -3019.125

$ws.Range("H116").Value = 1800
$ws.Range("I116").Value = 1666.6666
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 1666.6666
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 627.3334
$ws.Range("N116").Value = -6788

$ws.Range("H122").Value = 7288.609
$ws.Range("I122").Value = 8696.556
$ws.Range("K122").Value = 26089.668
$ws.Range("M122").Value = -23639.668

$ws.Range("H132").Value = 19234384
$ws.Range("I132").Value = 22730636
$ws.Range("K132").Value = 68191908
$ws.Range("M132").Value = -68189378

$ws.Range("H136").Value = 3181.88
$ws.Range("I136").Value = 3231.125
$ws.Range("K136").Value = 9693.375
$ws.Range("M136").Value = -7143.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1800
$ws.Range("I3").Value = 1666.6666
$ws.Range("J3").Value = 2200
$ws.Range("K3").Value = 1666.6666
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = -1552.6666
$ws.Range("N3").Value = -2428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 28980
$ws.Range("J63").Value = 28980
$ws.Range("L63").Value = 28980
$ws.Range("N63").Value = -30352

$ws.Range("H66").Value = 28980
$ws.Range("J66").Value = 28980
$ws.Range("L66").Value = 86940
$ws.Range("N66").Value = -93804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 150
$ws.Range("I17").Value = 100
$ws.Range("K17").Value = 300
$ws.Range("M17").Value = -131

$ws.Range("H34").Value = 746
$ws.Range("I34").Value = 102
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 306
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -222
$ws.Range("N34").Value = -9168

$ws.Range("H55").Value = 442.6154
$ws.Range("I55").Value = 151.14285
$ws.Range("J55").Value = 550
$ws.Range("K55").Value = 453.42855
$ws.Range("L55").Value = 1650
$ws.Range("M55").Value = -276.42855
$ws.Range("N55").Value = -2004

$ws.Range("H68").Value = 1058.5555
$ws.Range("I68").Value = 1200
$ws.Range("J68").Value = 968.5454999999999
$ws.Range("K68").Value = 3600
$ws.Range("L68").Value = 2905.6365
$ws.Range("M68").Value = -2789
$ws.Range("N68").Value = -4527.6365

$ws.Range("H71").Value = 1058.5555
$ws.Range("I71").Value = 1200
$ws.Range("J71").Value = 968.5454999999999
$ws.Range("K71").Value = 10800
$ws.Range("L71").Value = 8716.9095
$ws.Range("M71").Value = -6744
$ws.Range("N71").Value = -16828.9095

$ws.Range("H132").Value = 1430.875
$ws.Range("I132").Value = 687.5
$ws.Range("J132").Value = 1802.5625
$ws.Range("K132").Value = 6187.5
$ws.Range("L132").Value = 16223.0625
$ws.Range("M132").Value = -3657.5
$ws.Range("N132").Value = -21283.0625

$ws.Range("H137").Value = 4225.9614
$ws.Range("I137").Value = 2981.5386
$ws.Range("J137").Value = 5470.385
$ws.Range("K137").Value = 8944.6158
$ws.Range("L137").Value = 16411.155
$ws.Range("M137").Value = -3844.6158
$ws.Range("N137").Value = -26611.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3832.6365
$ws.Range("I102").Value = 4927.143
$ws.Range("J102").Value = 1917.25
$ws.Range("K102").Value = 4927.143
$ws.Range("L102").Value = 1917.25
$ws.Range("M102").Value = -3305.143
$ws.Range("N102").Value = -5161.25

$ws.Range("H122").Value = 2167.0476
$ws.Range("I122").Value = 2162.1538
$ws.Range("J122").Value = 2175
$ws.Range("K122").Value = 6486.4614
$ws.Range("L122").Value = 6525
$ws.Range("M122").Value = -4036.4614
$ws.Range("N122").Value = -11425

$ws.Range("H126").Value = 3215.7368
$ws.Range("I126").Value = 1969.9
$ws.Range("K126").Value = 5909.700000000001
$ws.Range("M126").Value = -3439.700000000001

$ws.Range("H133").Value = 77996.336
$ws.Range("J133").Value = 77996.336
$ws.Range("L133").Value = 77996.336
$ws.Range("N133").Value = -88116.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7080
$ws.Range("I40").Value = 17199.666
$ws.Range("J40").Value = 2743
$ws.Range("K40").Value = 17199.666
$ws.Range("L40").Value = 2743
$ws.Range("M40").Value = -17063.666
$ws.Range("N40").Value = -3015

$ws.Range("H122").Value = 9906.25
$ws.Range("I122").Value = 12560
$ws.Range("K122").Value = 37680
$ws.Range("M122").Value = -35230

$ws.Range("H136").Value = 3006.0923
$ws.Range("I136").Value = 1451.569
$ws.Range("J136").Value = 15886.429
$ws.Range("K136").Value = 4354.707
$ws.Range("L136").Value = 47659.287
$ws.Range("N136").Value = -52759.287
